$wb = $excel.ActiveWorkbook

# New trade row to append (row 30) on both the "All Trades" and "MarketMaking" sheets.
$sheetNames = @("All Trades", "MarketMaking")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $row = 30

    $ws.Cells.Item($row, 1).Value = 29
    # Use a leading apostrophe so Excel stores the date-looking string as literal
    # text instead of auto-converting it to a date serial value.
    $ws.Cells.Item($row, 2).Value = "'2026-02-17"
    $ws.Cells.Item($row, 3).Value = "23:58:07"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.53
    $ws.Cells.Item($row, 7).Value = ""
    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 101.1036569789373
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = ""
    $ws.Cells.Item($row, 17).Value = 0
}
